# Apply the "add delete classroom when import file" edit:
# - Update class ("Kelas") values for existing students from "1c" to "1e"
# - Add two new student rows (Eliminator / Orico) with class "2a"
# - Move the active selection to B4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: Eliminator (class "2a" typed first, matches new-string intern order)
$ws.Range("B4").Value = "2a"
$ws.Range("A4").Value = "Eliminator "
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 100

# New row 5: Orico
$ws.Range("A5").Value = "Orico"
$ws.Range("B5").Value = "2a"
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 70
$ws.Range("E5").Value = 80
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = 100

# Update existing rows: Kelas column (B) changes from "1c" to "1e"
$ws.Range("B2").Value = "1e"
$ws.Range("B3").Value = "1e"

# Update selection to B4
$ws.Range("B4").Select()
